$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data collection (e.g. a weather/garden log) picked up
# seven more days of observations (rows 401-407) below the prior last
# row (400). Column A carries the same date-formatted style as the rows
# above it, so copy that number format down first, then fill in values.
$ws.Cells.Item(400, 1).Copy()
$ws.Range("A401:A407").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 401
$ws.Cells.Item(401, 1).Value = 45844
$ws.Cells.Item(401, 2).Value = "Flowering"
$ws.Cells.Item(401, 3).Value = "Large"
$ws.Cells.Item(401, 4).Value = 70
$ws.Cells.Item(401, 5).Value = 91
$ws.Cells.Item(401, 6).Formula = "=ABS(D401-E401)"
$ws.Cells.Item(401, 7).Value = 0
$ws.Cells.Item(401, 8).Value = 0.1
$ws.Cells.Item(401, 9).Value = "No"
$ws.Cells.Item(401, 10).Value = 2
$ws.Cells.Item(401, 11).Value = "Bright"
$ws.Cells.Item(401, 12).Value = 9
$ws.Cells.Item(401, 13).Value = 0.44
$ws.Cells.Item(401, 14).Value = 66
$ws.Cells.Item(401, 15).Value = 29.98
$ws.Cells.Item(401, 16).Value = 10
$ws.Cells.Item(401, 17).Value = 0.1
$ws.Cells.Item(401, 18).Value = 9.9
$ws.Cells.Item(401, 19).Value = 62
$ws.Cells.Item(401, 20).Value = 0

# Row 402
$ws.Cells.Item(402, 1).Value = 45844
$ws.Cells.Item(402, 2).Value = "Nonflowering"
$ws.Cells.Item(402, 3).Value = "Medium"
$ws.Cells.Item(402, 4).Value = 70
$ws.Cells.Item(402, 5).Value = 91
$ws.Cells.Item(402, 6).Formula = "=ABS(D402-E402)"
$ws.Cells.Item(402, 7).Value = 0
$ws.Cells.Item(402, 8).Value = 0.1
$ws.Cells.Item(402, 9).Value = "No"
$ws.Cells.Item(402, 10).Value = 3
$ws.Cells.Item(402, 11).Value = "Bright"
$ws.Cells.Item(402, 12).Value = 9
$ws.Cells.Item(402, 13).Value = 0.44
$ws.Cells.Item(402, 14).Value = 66
$ws.Cells.Item(402, 15).Value = 29.98
$ws.Cells.Item(402, 16).Value = 10
$ws.Cells.Item(402, 17).Value = 0.1
$ws.Cells.Item(402, 18).Value = 9.9
$ws.Cells.Item(402, 19).Value = 62
$ws.Cells.Item(402, 20).Value = 0

# Row 403
$ws.Cells.Item(403, 1).Value = 45844
$ws.Cells.Item(403, 2).Value = "Nonflowering"
$ws.Cells.Item(403, 3).Value = "Small"
$ws.Cells.Item(403, 4).Value = 70
$ws.Cells.Item(403, 5).Value = 91
$ws.Cells.Item(403, 6).Formula = "=ABS(D403-E403)"
$ws.Cells.Item(403, 7).Value = 0
$ws.Cells.Item(403, 8).Value = 0.2
$ws.Cells.Item(403, 9).Value = "No"
$ws.Cells.Item(403, 10).Value = 3
$ws.Cells.Item(403, 11).Value = "Neutral"
$ws.Cells.Item(403, 12).Value = 9
$ws.Cells.Item(403, 13).Value = 0.44
$ws.Cells.Item(403, 14).Value = 66
$ws.Cells.Item(403, 15).Value = 29.98
$ws.Cells.Item(403, 16).Value = 10
$ws.Cells.Item(403, 17).Value = 0.1
$ws.Cells.Item(403, 18).Value = 9.9
$ws.Cells.Item(403, 19).Value = 62
$ws.Cells.Item(403, 20).Value = 0

# Row 404
$ws.Cells.Item(404, 1).Value = 45844
$ws.Cells.Item(404, 2).Value = "Nonflowering"
$ws.Cells.Item(404, 3).Value = "Medium"
$ws.Cells.Item(404, 4).Value = 70
$ws.Cells.Item(404, 5).Value = 91
$ws.Cells.Item(404, 6).Formula = "=ABS(D404-E404)"
$ws.Cells.Item(404, 7).Value = 0
$ws.Cells.Item(404, 8).Value = 0
$ws.Cells.Item(404, 9).Value = "No"
$ws.Cells.Item(404, 10).Value = 3
$ws.Cells.Item(404, 11).Value = "Neutral"
$ws.Cells.Item(404, 12).Value = 9
$ws.Cells.Item(404, 13).Value = 0.44
$ws.Cells.Item(404, 14).Value = 66
$ws.Cells.Item(404, 15).Value = 29.98
$ws.Cells.Item(404, 16).Value = 10
$ws.Cells.Item(404, 17).Value = 0.1
$ws.Cells.Item(404, 18).Value = 9.9
$ws.Cells.Item(404, 19).Value = 62
$ws.Cells.Item(404, 20).Value = 0

# Row 405
$ws.Cells.Item(405, 1).Value = 45844
$ws.Cells.Item(405, 2).Value = "Nonflowering"
$ws.Cells.Item(405, 3).Value = "Medium"
$ws.Cells.Item(405, 4).Value = 70
$ws.Cells.Item(405, 5).Value = 91
$ws.Cells.Item(405, 6).Formula = "=ABS(D405-E405)"
$ws.Cells.Item(405, 7).Value = 0
$ws.Cells.Item(405, 8).Value = 0
$ws.Cells.Item(405, 9).Value = "No"
$ws.Cells.Item(405, 10).Value = 3
$ws.Cells.Item(405, 11).Value = "Bright"
$ws.Cells.Item(405, 12).Value = 9
$ws.Cells.Item(405, 13).Value = 0.44
$ws.Cells.Item(405, 14).Value = 66
$ws.Cells.Item(405, 15).Value = 29.98
$ws.Cells.Item(405, 16).Value = 10
$ws.Cells.Item(405, 17).Value = 0.1
$ws.Cells.Item(405, 18).Value = 9.9
$ws.Cells.Item(405, 19).Value = 62
$ws.Cells.Item(405, 20).Value = 0

# Row 406
$ws.Cells.Item(406, 1).Value = 45844
$ws.Cells.Item(406, 2).Value = "Nonflowering"
$ws.Cells.Item(406, 3).Value = "Large"
$ws.Cells.Item(406, 4).Value = 70
$ws.Cells.Item(406, 5).Value = 91
$ws.Cells.Item(406, 6).Formula = "=ABS(D406-E406)"
$ws.Cells.Item(406, 7).Value = 0
$ws.Cells.Item(406, 8).Value = 0.15
$ws.Cells.Item(406, 9).Value = "No"
$ws.Cells.Item(406, 10).Value = 4
$ws.Cells.Item(406, 11).Value = "Neutral"
$ws.Cells.Item(406, 12).Value = 9
$ws.Cells.Item(406, 13).Value = 0.44
$ws.Cells.Item(406, 14).Value = 66
$ws.Cells.Item(406, 15).Value = 29.98
$ws.Cells.Item(406, 16).Value = 10
$ws.Cells.Item(406, 17).Value = 0.1
$ws.Cells.Item(406, 18).Value = 9.9
$ws.Cells.Item(406, 19).Value = 62
$ws.Cells.Item(406, 20).Value = 0

# Row 407
$ws.Cells.Item(407, 1).Value = 45844
$ws.Cells.Item(407, 2).Value = "Tree"
$ws.Cells.Item(407, 3).Value = "Medium"
$ws.Cells.Item(407, 4).Value = 70
$ws.Cells.Item(407, 5).Value = 91
$ws.Cells.Item(407, 6).Formula = "=ABS(D407-E407)"
$ws.Cells.Item(407, 7).Value = 0
$ws.Cells.Item(407, 8).Value = 0
$ws.Cells.Item(407, 9).Value = "No"
$ws.Cells.Item(407, 10).Value = 1
$ws.Cells.Item(407, 11).Value = "Bright"
$ws.Cells.Item(407, 12).Value = 9
$ws.Cells.Item(407, 13).Value = 0.44
$ws.Cells.Item(407, 14).Value = 66
$ws.Cells.Item(407, 15).Value = 29.98
$ws.Cells.Item(407, 16).Value = 10
$ws.Cells.Item(407, 17).Value = 0.1
$ws.Cells.Item(407, 18).Value = 9.9
$ws.Cells.Item(407, 19).Value = 62
$ws.Cells.Item(407, 20).Value = 0

# Reflect the scroll/selection state Excel saved after entering this
# block: the view is scrolled so row 398 is at the top and the active
# cell sits just past the last populated row, in the last column.
$excel.ActiveWindow.ScrollRow = 398
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("U401").Select()

